$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P5").Value = '" ~OT ~ = Petron Live Implementation EFG Marketing ~OB Others|Petron Live Implementation - EFG Marketing San Pablo Laguna May 8, 9, and 11.| R "'
$ws.Range("P6").Value = '" ~OB Others|Petron Live Implementation - EFG Marketing San Pablo Laguna May 8, 9, and 11.| R "'
$ws.Range("P7").Value = '" ~OT ~ = Petron Live Implementation - San Pablo Laguna ~OB Others|Petron Live Implementation - EFG Marketing San Pablo Laguna May 8, 9, and 11.| R "'
$ws.Range("P8").Value = '" ~OB Others|Petron Live Implementation - BMF Gasul May 13, 14 and 15.  May 12 - Transportation from Manila to Dagupan| R "'
$ws.Range("P9").Value = '" ~OB Others|Petron Live Implementation - BMF Gasul May 13, 14 and 15.  May 12 - Transportation from Manila to Dagupan| R "'
$ws.Range("P10").Value = '" ~OB Others|Petron Live Implementation - BMF Gasul May 13, 14 and 15.  May 12 - Transportation from Manila to Dagupan| R "'
$ws.Range("P11").Value = '" ~OT ~ = Petron Live Implementation - BMF Gasul, Dagupan City ~OB Others|Petron Live Implementation - BMF Gasul May 13, 14 and 15.  May 12 - Transportation from Manila to Dagupan| R "'
$ws.Range("P12").Value = '" ~OT ~ = Petron Live Implementation - BMF Gasul, Dagupan City"'
